$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.224.95"
$ws.Range("E2").Value = "  +5.71%  "
$ws.Range("D3").Value = "2.657.36"
$ws.Range("E3").Value = "  +7.97%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'185.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.33%  "
$ws.Range("D6").Value = "'583.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("D9").Value = "'0.193"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.09%  "
$ws.Range("D10").Value = "2.655.90"
$ws.Range("E10").Value = "  +7.99%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +5.00%  "
$ws.Range("D13").Value = "'4.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "3.147.67"
$ws.Range("E14").Value = "  +8.28%  "
$ws.Range("D15").Value = "74.143.98"
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "'26.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.51%  "
$ws.Range("D18").Value = "2.664.15"
$ws.Range("E18").Value = "  +8.54%  "
$ws.Range("E19").Value = "  +28.93%  "
$ws.Range("D20").Value = "'11.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.90%  "
$ws.Range("D21").Value = "'370.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.56%  "
$ws.Range("D22").Value = "'2.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.82%  "
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").Value = "'6.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'69.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("D28").Value = "'9.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.71%  "
$ws.Range("D29").Value = "2.795.67"
$ws.Range("E29").Value = "  +8.27%  "
$ws.Range("E30").Value = "  -19.89%  "
$ws.Range("D31").Value = "0.0₃0934"
$ws.Range("E31").Value = "  +6.97%  "
$ws.Range("D32").Value = "'516.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.31%  "
$ws.Range("E33").Value = "  +11.09%  "
$ws.Range("D34").Value = "'7.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D37").Value = "'163.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("E38").Value = "  +5.71%  "
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'170.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +27.54%  "
$ws.Range("D43").Value = "'4.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.79%  "
$ws.Range("D44").Value = "'0.326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.44%  "
$ws.Range("E45").Value = "  +6.44%  "
$ws.Range("E46").Value = "  +7.69%  "
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").Value = "'2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.98%  "
$ws.Range("D49").Value = "'0.0841"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.70%  "
$ws.Range("D50").Value = "'3.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("E51").Value = "  +6.18%  "
